$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update tag values that now contain multiple comma-separated tags ---
# (Order matters: it determines the order new shared strings are appended in.)

# Row 96: "Describe what Human Resource means to you." -> tags gain "behavioral"
$ws.Range("B96").Value = "general, behavioral"

# Row 2: "Tell me about yourself." -> tags become general, personal / amazon, microsoft
$ws.Range("C2").Value = "amazon, microsoft"
$ws.Range("B2").Value = "general, personal"

# --- Update the Excel->JSON formula so it SUBSTITUTEs ", " with "", "" inside
#     the tag columns before re-quoting them, letting a single cell hold
#     several comma-separated tags. ---

# Set the second block FIRST so it claims shared-formula id si=1 and leaves
# id si=0 free for the pre-existing E3:E66 group re-created just below
# (matches the si numbering produced by the real edit).

# Second formula block (rows 67..117, new shared group si=1 anchored at E67)
$ws.Range("E67:E117").Formula = '=CONCATENATE("{ question : """,A67,""", tags : [""",SUBSTITUTE(B67, ", ", """, """),"""",IF(C67="null","",_xlfn.CONCAT(", """,SUBSTITUTE(C67, ", ", """, """),"""" )), "] },")'

# First formula block (rows 2 and 3..66, shared group si=0 anchored at E3)
$ws.Range("E2").Formula = '=CONCATENATE("{ question : """,A2,""", tags : [""",SUBSTITUTE(B2, ", ", """, """),"""",IF(C2="null","",_xlfn.CONCAT(", """,SUBSTITUTE(C2, ", ", """, """),"""" )), "] },")'
$ws.Range("E3:E66").Formula = '=CONCATENATE("{ question : """,A3,""", tags : [""",SUBSTITUTE(B3, ", ", """, """),"""",IF(C3="null","",_xlfn.CONCAT(", """,SUBSTITUTE(C3, ", ", """, """),"""" )), "] },")'

# --- Column sizing for the newly-meaningful tag columns B and C ---
$ws.Columns.Item(2).ColumnWidth = 15.8333333333
$ws.Columns.Item(3).ColumnWidth = 14.3333333333

# --- Selection / scroll state ---
$ws.Range("E3").Select() | Out-Null
